# Auto-generated Excel COM-interop script to apply value updates
# per the Pandaemonium_Profits diff (scheduled-runner price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1476.625
$ws.Range("I11").Value = 1476.625
$ws.Range("K11").Value = 1476.625
$ws.Range("M11").Value = -1336.625
$ws.Range("H33").Value = 207.85294
$ws.Range("I33").Value = 183.39285
$ws.Range("K33").Value = 183.39285
$ws.Range("M33").Value = 45.60714999999999
$ws.Range("H107").Value = 800.3913
$ws.Range("I107").Value = 663.3684
$ws.Range("J107").Value = 1451.25
$ws.Range("K107").Value = 663.3684
$ws.Range("L107").Value = 1451.25
$ws.Range("M107").Value = 1256.6316
$ws.Range("N107").Value = -5291.25
$ws.Range("H132").Value = 1259.1072
$ws.Range("I132").Value = 1211.6154
$ws.Range("K132").Value = 3634.8462
$ws.Range("M132").Value = -1104.8462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5517.25
$ws.Range("I61").Value = 4422.8
$ws.Range("J61").Value = 8253.375
$ws.Range("K61").Value = 4422.8
$ws.Range("L61").Value = 8253.375
$ws.Range("M61").Value = -4210.8
$ws.Range("N61").Value = -8677.375
$ws.Range("H63").Value = 3022
$ws.Range("I63").Value = 3022
$ws.Range("K63").Value = 3022
$ws.Range("M63").Value = -2336
$ws.Range("H66").Value = 3022
$ws.Range("I66").Value = 3022
$ws.Range("K66").Value = 15110
$ws.Range("M66").Value = -11678
$ws.Range("H122").Value = 1713.2858
$ws.Range("I122").Value = 1738.3182
$ws.Range("J122").Value = 1621.5
$ws.Range("K122").Value = 5214.9546
$ws.Range("L122").Value = 4864.5
$ws.Range("M122").Value = -2764.9546
$ws.Range("N122").Value = -9764.5
$ws.Range("H136").Value = 5517.25
$ws.Range("I136").Value = 4422.8
$ws.Range("J136").Value = 8253.375
$ws.Range("K136").Value = 13268.4
$ws.Range("L136").Value = 24760.125
$ws.Range("M136").Value = -10718.4
$ws.Range("N136").Value = -29860.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 212
$ws.Range("I22").Value = 212
$ws.Range("K22").Value = 212
$ws.Range("M22").Value = -39
$ws.Range("H86").Value = 2926.5
$ws.Range("I86").Value = 3130.2856
$ws.Range("K86").Value = 3130.2856
$ws.Range("M86").Value = -2007.2856
$ws.Range("H89").Value = 2926.5
$ws.Range("I89").Value = 3130.2856
$ws.Range("K89").Value = 15651.428
$ws.Range("M89").Value = -10035.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2520.5
$ws.Range("I31").Value = 1714.64
$ws.Range("K31").Value = 1714.64
$ws.Range("M31").Value = -1419.64
$ws.Range("H34").Value = 2520.5
$ws.Range("I34").Value = 1714.64
$ws.Range("K34").Value = 1714.64
$ws.Range("M34").Value = -1512.64
$ws.Range("H69").Value = 19450
$ws.Range("I69").Value = 15933.333
$ws.Range("J69").Value = 30000
$ws.Range("K69").Value = 15933.333
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = -15184.333
$ws.Range("N69").Value = -31498
$ws.Range("H72").Value = 19450
$ws.Range("I72").Value = 15933.333
$ws.Range("J72").Value = 30000
$ws.Range("K72").Value = 47799.999
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = -44055.999
$ws.Range("N72").Value = -97488
$ws.Range("H99").Value = 1485.9412
$ws.Range("I99").Value = 1563.5
$ws.Range("J99").Value = 1299.8
$ws.Range("K99").Value = 1563.5
$ws.Range("L99").Value = 1299.8
$ws.Range("M99").Value = -65.5
$ws.Range("N99").Value = -4295.8
$ws.Range("H105").Value = 1200.091
$ws.Range("H126").Value = 1485.9412
$ws.Range("I126").Value = 1563.5
$ws.Range("J126").Value = 1299.8
$ws.Range("K126").Value = 4690.5
$ws.Range("L126").Value = 3899.4
$ws.Range("M126").Value = -2220.5
$ws.Range("N126").Value = -8839.4
$ws.Range("H132").Value = 2471.9092
$ws.Range("I132").Value = 2091.4707
$ws.Range("J132").Value = 3765.4
$ws.Range("K132").Value = 6274.4121
$ws.Range("L132").Value = 11296.2
$ws.Range("M132").Value = -3744.4121
$ws.Range("N132").Value = -16356.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 789.4
$ws.Range("I5").Value = 519.25
$ws.Range("K5").Value = 1557.75
$ws.Range("M5").Value = -1445.75
$ws.Range("H113").Value = 735.41
$ws.Range("I113").Value = 754.1786
$ws.Range("J113").Value = 636.875
$ws.Range("K113").Value = 2262.5358
$ws.Range("L113").Value = 1910.625
$ws.Range("M113").Value = -92.53579999999965
$ws.Range("N113").Value = -6250.625
$ws.Range("H122").Value = 996.6
$ws.Range("I122").Value = 550.6667
$ws.Range("J122").Value = 1293.8889
$ws.Range("K122").Value = 4956.0003
$ws.Range("L122").Value = 11645.0001
$ws.Range("M122").Value = -2506.0003
$ws.Range("N122").Value = -16545.0001
$ws.Range("H135").Value = 789.4
$ws.Range("I135").Value = 519.25
$ws.Range("K135").Value = 4673.25
$ws.Range("M135").Value = -2138.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2475.3125
$ws.Range("I126").Value = 1720.6666
$ws.Range("J126").Value = 3141.1765
$ws.Range("K126").Value = 5161.9998
$ws.Range("L126").Value = 9423.529500000001
$ws.Range("M126").Value = -2691.9998
$ws.Range("N126").Value = -14363.5295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 9800
$ws.Range("J14").Value = 9800
$ws.Range("L14").Value = 9800
$ws.Range("N14").Value = -10144
$ws.Range("H22").Value = 226
$ws.Range("I22").Value = 200.07692
$ws.Range("J22").Value = 259.7
$ws.Range("K22").Value = 200.07692
$ws.Range("L22").Value = 259.7
$ws.Range("M22").Value = 94.92308
$ws.Range("N22").Value = -849.7
$ws.Range("H27").Value = 226
$ws.Range("I27").Value = 200.07692
$ws.Range("J27").Value = 259.7
$ws.Range("K27").Value = 200.07692
$ws.Range("L27").Value = 259.7
$ws.Range("M27").Value = -93.07692
$ws.Range("N27").Value = -473.7
$ws.Range("H132").Value = 5376.2
$ws.Range("I132").Value = 5157.875
$ws.Range("J132").Value = 6249.5
$ws.Range("K132").Value = 15473.625
$ws.Range("L132").Value = 18748.5
$ws.Range("M132").Value = -12943.625
$ws.Range("N132").Value = -23808.5
$ws.Range("H136").Value = 4974.1763
$ws.Range("I136").Value = 2985.9583
$ws.Range("J136").Value = 9745.9
$ws.Range("K136").Value = 8957.874899999999
$ws.Range("L136").Value = 29237.7
$ws.Range("M136").Value = -6407.874899999999
$ws.Range("N136").Value = -34337.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2239.0645
$ws.Range("I122").Value = 1936.24
$ws.Range("J122").Value = 3500.8333
$ws.Range("K122").Value = 5808.72
$ws.Range("L122").Value = 10502.4999
$ws.Range("M122").Value = -3407.3044
$ws.Range("N122").Value = -15402.4999
